$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.253.93"
$ws.Range("E2").Value = "  +0.62%  "
# Row 3
$ws.Range("D3").Value = "3.184.83"
$ws.Range("E3").Value = "  -0.11%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.04%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.75%  "
# Row 7
$ws.Range("E7").Value = "  -0.09%  "
# Row 8
$ws.Range("D8").Value = "3.184.83"
$ws.Range("E8").Value = "  -0.08%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.25%  "
# Row 10
$ws.Range("E10").Value = "  -0.41%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.02%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.508"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.61%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "
# Row 15
$ws.Range("D15").Value = "3.704.02"
$ws.Range("E15").Value = "  -0.18%  "
# Row 16
$ws.Range("D16").Value = "66.214.79"
$ws.Range("E16").Value = "  +0.48%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "
# Row 18
$ws.Range("D18").Value = "3.183.30"
$ws.Range("E18").Value = "  -0.39%  "
# Row 19
$ws.Range("E19").Value = "  +1.15%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.02%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.10%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.91%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.19%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "
# Row 26
$ws.Range("E26").Value = "  +0.06%  "
# Row 27
$ws.Range("E27").Value = "  +0.43%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.50%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.23%  "
# Row 30
$ws.Range("E30").Value = "  +5.97%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.82%  "
# Row 33
$ws.Range("E33").Value = "  +0.02%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.64%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.65%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.39%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0899"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.14%  "
# Row 39
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0723"
$ws.Range("E39").Value = "  +11.64%  "
# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.38%  "
# Row 41
$ws.Range("E41").Value = "  +5.20%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.13%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.77%  "
# Row 45
$ws.Range("E45").Value = "  +2.46%  "
# Row 46
$ws.Range("D46").Value = "2.832.82"
$ws.Range("E46").Value = "  -3.36%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "
# Row 48
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.94%  "
# Row 49
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.117"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.63%  "
